$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 169 (this pushes the existing
# rows 169..306 down to 170..307, growing the sheet's dimension from
# A1:R306 to A1:R307).
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with its values.
$ws.Cells.Item(169, 1).Value  = 5
$ws.Cells.Item(169, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(169, 3).Value  = "Maule"
$ws.Cells.Item(169, 4).Value  = 44673
$ws.Cells.Item(169, 5).Value  = 7
$ws.Cells.Item(169, 6).Value  = 100114014
$ws.Cells.Item(169, 7).Value  = "Betarraga"
$ws.Cells.Item(169, 8).Value  = "Sin especificar"
$ws.Cells.Item(169, 9).Value  = "Primera"
$ws.Cells.Item(169, 10).Value = 5000
$ws.Cells.Item(169, 11).Value = 700
$ws.Cells.Item(169, 12).Value = 700
$ws.Cells.Item(169, 13).Value = 700
$ws.Cells.Item(169, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(169, 15).Value = "Región del Maule"
$ws.Cells.Item(169, 16).Value = 140
$ws.Cells.Item(169, 17).Value = 5
$ws.Cells.Item(169, 18).Value = "Hortaliza"
